$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 12 data rows (rows 2-13), shifting remaining data rows up.
$ws.Rows("2:13").Delete()

# New rows appended after the shifted data (rows 11-21 of the remaining 9 + 11 new = 20 data rows, rows 2-21)
$newData = @(
    @(1.691281750798222, -3.219833493232723, 15.28036725521086),
    @(0.7758507728576624, -2.317984580993662, 10.35012340545653),
    @(-0.05869728326797519, -4.645559132099153, 8.129511415958403),
    @(-7.922254800796521, -8.533369660377506, 4.239331245422359),
    @(-7.12406146526336, -9.578974485397339, 5.120484650135045),
    @(-4.634187221527098, -9.333477973937988, 4.373982667922966),
    @(-2.660379245877256, -9.284614562988283, -3.805972993373912),
    @(0.6678269803524017, -9.062278509140016, -10.71978342533112),
    @(2.381343364715557, -7.552583456039422, 0.6790638566017475),
    @(-2.60627746582033, -5.606535911560057, 9.005005836486802),
    @(-6.737108409404774, -5.243253648281096, 5.65115070343016)
)

$startRow = 11
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}
